# Adds example traversal orders / explanatory text to the "General Tree"
# section of the document (the three lines that originally read
# "Preorder: root, son, next", "Inorder: " and "Postorder:" with nothing,
# or little, after the colon).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: "Preorder: root, son, next" -> append explanatory note
# ---------------------------------------------------------------------
$rngPre = $d.Content
$foundPre = $rngPre.Find.Execute("Preorder: root, son, next", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundPre) {
    throw "Could not find 'Preorder: root, son, next'"
}
$rngPre.Collapse(0)
$rngPre.InsertAfter(" (very standard way of traversing, top to bottom, left to right)")
$afterPre = $rngPre.End

# ---------------------------------------------------------------------
# Hunk 2: "Inorder: " (General Tree, empty after the colon) -> "son, root, next"
# Re-use/extend the existing ": " run so the inserted text keeps the exact
# same run formatting (Times New Roman / sz 24 / szCs 24) as the diff.
# ---------------------------------------------------------------------
$rngInScope = $d.Range($afterPre, $d.Content.End)
$foundIn = $rngInScope.Find.Execute(": ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundIn) {
    throw "Could not find 'Inorder: ' after the Preorder line"
}
$rngInScope.Text = ": son, root, next"
$afterIn = $rngInScope.End

# ---------------------------------------------------------------------
# Hunk 3: "Postorder:" (General Tree, nothing after the colon) ->
#          " son, next, root (note this is tricky. study examples carefully)"
# Same technique: extend the existing ":" run in place.
# ---------------------------------------------------------------------
$rngPostScope = $d.Range($afterIn, $d.Content.End)
$foundPost = $rngPostScope.Find.Execute(":", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundPost) {
    throw "Could not find 'Postorder:' after the Inorder line"
}
$rngPostScope.Text = ": son, next, root (note this is tricky. study examples carefully)"

Write-Host "Done applying General Tree traversal annotations"
